$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '86.517.06'
$ws.Range("E2").Value = '  +5.19%  '

# Row 3
$ws.Range("D3").Value = '3.253.57'
$ws.Range("E3").Value = '  +2.42%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '209.34'
$ws.Range("E5").Value = '  -3.29%  '

# Row 6
$ws.Range("D6").Value = '623.42'
$ws.Range("E6").Value = '  +0.34%  '

# Row 7
$ws.Range("D7").Value = '0.372'
$ws.Range("E7").Value = '  +28.15%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").Value = '0.645'
$ws.Range("E9").Value = '  +10.52%  '

# Row 10
$ws.Range("D10").Value = '3.245.19'
$ws.Range("E10").Value = '  +2.12%  '

# Row 11
$ws.Range("D11").Value = '0.572'
$ws.Range("E11").Value = '  -7.65%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.177'
$ws.Range("E12").Value = '  +7.36%  '

# Row 13
$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").Value = '0.0000257'
$ws.Range("E13").Value = '  -0.31%  '

# Row 14
$ws.Range("D14").Value = '33.93'
$ws.Range("E14").Value = '  +4.85%  '

# Row 15
$ws.Range("D15").Value = '3.838.01'
$ws.Range("E15").Value = '  +1.92%  '

# Row 16
$ws.Range("D16").Value = '5.22'
$ws.Range("E16").Value = '  -1.40%  '

# Row 17
$ws.Range("D17").Value = '86.402.19'
$ws.Range("E17").Value = '  +5.37%  '

# Row 18
$ws.Range("D18").Value = '3.238.68'
$ws.Range("E18").Value = '  +2.47%  '

# Row 19
$ws.Range("D19").Value = '13.98'
$ws.Range("E19").Value = '  -0.44%  '

# Row 20
$ws.Range("D20").Value = '2.95'
$ws.Range("E20").Value = '  -8.13%  '

# Row 21
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '428.11'
$ws.Range("E21").Value = '  -1.58%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '8.95'
$ws.Range("E22").Value = '  +0.05%  '

# Row 23
$ws.Range("D23").Value = '5.35'
$ws.Range("E23").Value = '  +4.09%  '

# Row 24
$ws.Range("E24").Value = '  -1.72%  '

# Row 25
$ws.Range("D25").Value = '12.15'
$ws.Range("E25").Value = '  +4.95%  '

# Row 26
$ws.Range("D26").Value = '5.11'
$ws.Range("E26").Value = '  -4.13%  '

# Row 27
$ws.Range("D27").Value = '3.430.00'
$ws.Range("E27").Value = '  +2.56%  '

# Row 28
$ws.Range("D28").Value = '75.91'
$ws.Range("E28").Value = '  -1.15%  '

# Row 29
$ws.Range("D29").Value = '0.0000128'
$ws.Range("E29").Value = '  +4.98%  '

# Row 30
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = '  +0.25%  '

# Row 31
$ws.Range("D31").Value = '0.173'
$ws.Range("E31").Value = '  +12.13%  '

# Row 32
$ws.Range("E32").Value = '  +0.30%  '

# Row 33
$ws.Range("D33").Value = '8.83'
$ws.Range("E33").Value = '  -2.34%  '

# Row 34
$ws.Range("D34").Value = '545.46'
$ws.Range("E34").Value = '  -5.15%  '

# Row 35
$ws.Range("D35").Value = '1.42'
$ws.Range("E35").Value = '  -5.51%  '

# Row 36
$ws.Range("D36").Value = '1.95'
$ws.Range("E36").Value = '  -2.00%  '

# Row 37
$ws.Range("D37").Value = '6.73'
$ws.Range("E37").Value = '  +9.39%  '

# Row 38
$ws.Range("D38").Value = '0.136'
$ws.Range("E38").Value = '  -11.11%  '

# Row 39
$ws.Range("D39").Value = '22.44'
$ws.Range("E39").Value = '  -1.47%  '

# Row 40
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("D41").Value = '21.69'
$ws.Range("E41").Value = '  +3.95%  '

# Row 42
$ws.Range("D42").Value = '0.392'
$ws.Range("E42").Value = '  -3.51%  '

# Row 43
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").Value = '  -1.06%  '

# Row 44
$ws.Range("E44").Value = '  -5.44%  '

# Row 45
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.30%  '

# Row 46
$ws.Range("D46").Value = '155.19'
$ws.Range("E46").Value = '  -3.60%  '

# Row 47
$ws.Range("D47").Value = '177.19'
$ws.Range("E47").Value = '  -5.74%  '

# Row 48
$ws.Range("D48").Value = '1.32'
$ws.Range("E48").Value = '  -0.67%  '

# Row 49
$ws.Range("D49").Value = '44.06'
$ws.Range("E49").Value = '  -1.38%  '

# Row 50
$ws.Range("D50").Value = '4.24'
$ws.Range("E50").Value = '  +0.75%  '

# Row 51
$ws.Range("D51").Value = '0.622'
$ws.Range("E51").Value = '  -1.28%  '
